$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.621.43"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.889.15"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4900"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2938"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06699"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "1.890.29"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.134"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6678"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "30.569.61"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007858"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "2.165.38"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.314"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.21%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "190.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.199"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.503"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.929"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.394"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05237"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7437"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.098"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01822"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.685"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.065"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +31.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4419"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.916"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9932"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1376"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.547"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.038"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05839"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.422"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
